$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-15 moves forward one day:
# serial date 45202 (2023-10-03) -> 45203 (2023-10-04)
for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value2 = 45203
    }
}
